$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.979.60'
$ws.Range('E2').Value = '  +4.62%  '
$ws.Range('D3').Value = '1.914.46'
$ws.Range('E3').Value = '  +1.63%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.81'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.40%  '
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '47.86'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +10.81%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.376'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.41'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0761'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.39%  '
$ws.Range('E12').Value = '  +2.27%  '
$ws.Range('E13').Value = '  +13.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.823'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.27%  '
$ws.Range('D15').Value = '2.192.95'
$ws.Range('E15').Value = '  +1.54%  '
$ws.Range('E16').Value = '  +2.60%  '
$ws.Range('D17').Value = '1.916.35'
$ws.Range('E17').Value = '  +1.78%  '
$ws.Range('D18').Value = '37.012.14'
$ws.Range('E18').Value = '  +4.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '74.46'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.57%  '
$ws.Range('E20').Value = '  +3.91%  '
$ws.Range('E21').Value = '  +6.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '251.38'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.98%  '
$ws.Range('E23').Value = '  +0.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E25').Value = '  -7.03%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.53'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.35%  '
$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.15'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.80'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.70'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.52%  '
$ws.Range('E30').Value = '  +0.87%  '
$ws.Range('E31').Value = '  +6.48%  '
$ws.Range('E32').Value = '  +2.36%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0910'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +25.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.30'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.26%  '
$ws.Range('E35').Value = '  +1.13%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.26'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +40.64%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.888'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.57%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.47'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.11%  '
$ws.Range('E40').Value = '  +2.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '104.61'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.56%  '
$ws.Range('E42').Value = '  +3.95%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.53'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.69%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.86'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +18.86%  '
$ws.Range('E45').Value = '  +2.30%  '
$ws.Range('D46').Value = '1.350.50'
$ws.Range('E46').Value = '  +3.03%  '
$ws.Range('E47').Value = '  +1.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0833'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.18%  '
$ws.Range('E49').Value = '  +2.92%  '
$ws.Range('E50').Value = '  +1.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.78'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +13.53%  '
